$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-27 (A:H), reflecting the inserted/reordered rows
$data = @(
    @("nf_era_age", 0, 0, 0, 0.67, 1, 0.33, 0.33),
    @("nf_era_alcohol", 0.33, 0, 0, 0, 0.67, 0.67, 0.67),
    @("nf_era_cholesterol", 0, 0, 0, 0, 1, 1, 1),
    @("nf_era_gender", 0.67, 0, 0, 0.33, 0.33, 0, 0),
    @("nf_era_glucose", 0, 0, 0, 0.33, 1, 0.67, 0.67),
    @("nf_era_sg", 0, 0, 0, 1, 1, 0, 0),
    @("un_franzosa_ControlvsCD_Age", 0, 0, 0, 0.67, 1, 0.33, 0.33),
    @("un_franzosa_ControlvsCD_ConvCD", 0, 0, 0, 0.67, 1, 0.33, 0.33),
    @("un_franzosa_ControlvsCD_Fp", 0.33, 0, 0, 0.33, 0.67, 0.33, 0.33),
    @("un_franzosa_ControlvsDisease_Age", 0, 0, 0, 0.33, 1, 0.67, 0.67),
    @("un_franzosa_ControlvsDisease_ConvDisease", 0, 0, 0, 0.33, 1, 0.67, 0.67),
    @("un_franzosa_ControlvsDisease_Fp", 0.67, 0, 0, 0, 0.33, 0.33, 0.33),
    @("un_franzosa_ControlvsUC_Age", 0, 0, 0, 0.67, 1, 0.33, 0.33),
    @("un_franzosa_ControlvsUC_ConvUC", 0, 0, 0.33, 0.33, 0.67, 0.67, 0.33),
    @("un_franzosa_ControlvsUC_Fp", 0.33, 0, 0, 0.33, 0.67, 0.33, 0.33),
    @("nf_yachida_age", 0, 0, 0, 0.33, 1, 0.67, 0.67),
    @("nf_yachida_alcohol", 0, 0, 0, 0.67, 1, 0.33, 0.33),
    @("nf_yachida_gender", 0, 0, 0, 0.67, 1, 0.33, 0.33),
    @("nf_yachida_healthyvscancer", 0, 0, 0, 0.33, 1, 0.67, 0.67),
    @("nf_yachida_healthyvsstageIII_IV", 0, 0, 0, 0.67, 1, 0.33, 0.33),
    @("nf_wang_age", 0, 0, 0, 0.33, 1, 0.67, 0.67),
    @("nf_wang_bmi", 0, 0, 0, 0.67, 1, 0.33, 0.33),
    @("nf_wang_creatinine", 0.33, 0, 0, 0.67, 0.67, 0, 0),
    @("nf_wang_egfr", 0.33, 0, 0, 0.67, 0.67, 0, 0),
    @("nf_wang_studygroup", 0.33, 0, 0, 0.33, 0.67, 0.33, 0.33),
    @("nf_wang_urea", 0.33, 0, 0, 0.67, 0.67, 0, 0)
)

$numRows = $data.Count
$numCols = 8

# Build a 2D array (1-based) for bulk write
$arr = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}

$startRow = 2
$endRow = $startRow + $numRows - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $numCols))
$rng.Value = $arr

